# Refreshes the cryptocurrency price table ("Price" column D and
# "Volume(1h)" column E) on the active worksheet for rows 2-51, mirroring
# an automated GitHub Actions data refresh.
#
# Values are assigned as plain text (matching the workbook's existing
# inline-string cells). Some new Price values would otherwise be
# auto-detected by Excel as numbers (e.g. "192.45"), so for those cells
# we temporarily force a text number format before assigning the value
# and then restore the cell's style so formatting stays unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.874.90"
$ws.Range("E2").Value = "  +2.78%  "
$ws.Range("D3").Value = "3.376.77"
$ws.Range("E3").Value = "  +3.25%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "192.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "592.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.609"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.16%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.134"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.25%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.77"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.89%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.418"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.61%  "
$ws.Range("D12").Value = "3.965.48"
$ws.Range("E12").Value = "  +3.29%  "
$ws.Range("E13").Value = "  -1.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.66"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.02%  "
$ws.Range("D15").Value = "69.821.05"
$ws.Range("E15").Value = "  +2.71%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000172"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.56%  "
$ws.Range("D17").Value = "3.379.95"
$ws.Range("E17").Value = "  +2.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "456.72"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +16.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.84"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.42%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.79"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.51%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "76.25"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.53%  "
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").Value = "3.524.37"
$ws.Range("E24").Value = "  +3.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.524"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.29%  "
$ws.Range("E26").Value = "  +3.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.192"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.70%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.50"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.02"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "23.40"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.59"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.32%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.30"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.95%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.57"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "164.89"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.75%  "
$ws.Range("E38").Value = "  +2.66%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "27.44"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.50%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.813"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.24%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.61"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.42%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.58"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.61%  "
$ws.Range("D43").Value = "2.734.05"
$ws.Range("E43").Value = "  +3.93%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.55"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.66%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0692"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.64%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "25.55"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.41%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "41.02"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.71%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "337.48"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.36%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0285"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.84%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "32.77"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.02"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.24%  "
